# Update "Elapsed Duration(Hrs)" values (column G) across the R1-R6 sheets
# to reflect the latest refresh timestamps.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3929:44:03" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "69:16:41" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12111:07:44" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3240:51:13" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "479:02:47" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2956:57:33" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "184:09:48" },
    @{ Sheet = "R4"; Cell = "G4"; Value = "72:22:13" },
    @{ Sheet = "R4"; Cell = "G5"; Value = "69:59:46" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "430:56:32" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "71:28:50" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
